# Applies the commit's changes:
#  1. Bumps the publication Date on the Metadata sheet.
#  2. Swaps the two "Mapping" columns (AK <-> AL, i.e. columns 37 and 38) on
#     the Elements sheet: header text, data rows (2-6) and column widths.

$wb = $excel.ActiveWorkbook

# --- 1. Update the Date value on the Metadata sheet -------------------------
$metaWs = $wb.Worksheets.Item("Metadata")
$metaWs.Range("B8").Value = "2024-03-22T16:25:12+00:00"

# --- 2. Swap the "Mapping: RIM Mapping" (AK) and                              
#        "Mapping: Spécification métier vers l'extension ROR CommuneCog" (AL)
#        columns on the Elements sheet -------------------------------------
$elWs = $wb.Worksheets.Item("Elements")

# Only rows 1, 3, 5 and 6 actually differ between AK/AL (rows 2 and 4 are
# blank in both columns) - skip the no-op rows so we don't needlessly
# rewrite cells that are already correct.
$rowsToSwap = @(1, 3, 5, 6)

# Remember AK (col 37) / AL (col 38) values before overwriting.
$akValues = @{}
$alValues = @{}
foreach ($r in $rowsToSwap) {
    $akValues[$r] = $elWs.Cells.Item($r, 37).Value2
    $alValues[$r] = $elWs.Cells.Item($r, 38).Value2
}

foreach ($r in $rowsToSwap) {
    $elWs.Cells.Item($r, 37).Value = $alValues[$r]
    $elWs.Cells.Item($r, 38).Value = $akValues[$r]
}

# Swap the column widths too, so col 37 (AK) now has the (wider) width that
# used to belong to col 38/AL, and col 38 (AL) gets the narrower ex-AK width.
# (Original OOXML widths were 24.98046875 / 69.23828125; ColumnWidth only
# round-trips to the nearest 1/6-character snap grid, so these literals are
# chosen to land on the closest achievable grid points: ~25.0 / ~69.1667.)
$elWs.Columns.Item(37).ColumnWidth = 68.33
$elWs.Columns.Item(38).ColumnWidth = 24.165
